# Auto-generated edit script applying the cryptos.xlsx diff.
# For cells whose new text would otherwise be auto-parsed by Excel as a
# number (single decimal point, e.g. "316.01"), the cell is forced to
# Text format for the assignment and then restored to General, so the
# stored value stays a string exactly like the target inlineStr cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.793.96'
$ws.Range('E2').Value = '  +0.44%  '
$ws.Range('D3').Value = '2.294.86'
$ws.Range('E3').Value = '  -0.36%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.01'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  -0.35%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '104.16'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  +0.42%  '
$ws.Range('E7').Value = '  -1.11%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.602'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').Value = '  -1.43%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.45'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  -1.02%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0906'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  -0.73%  '
$ws.Range('E12').Value = '  +1.14%  '
$ws.Range('E13').Value = '  +2.45%  '
$ws.Range('E14').Value = '  +3.00%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.27'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  -0.52%  '
$ws.Range('D16').Value = '2.642.48'
$ws.Range('E16').Value = '  -0.37%  '
$ws.Range('D17').Value = '2.297.62'
$ws.Range('E17').Value = '  -0.07%  '
$ws.Range('D18').Value = '42.715.85'
$ws.Range('E18').Value = '  +0.46%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.72'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  +31.65%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.49'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  -0.40%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0000106'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  -0.33%  '
$ws.Range('E22').Value = '  +1.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.56'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  -0.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '264.10'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  -5.35%  '
$ws.Range('E25').Value = '  -3.13%  '
$ws.Range('E26').Value = '  +0.42%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.88'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  +0.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.35'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  -0.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.87'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  +16.69%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.36'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  -1.87%  '
$ws.Range('E31').Value = '  +3.64%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '166.76'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  +1.00%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0874'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  -0.42%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.131'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  -4.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.57'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  -0.91%  '
$ws.Range('E36').Value = '  -3.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.58'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  -0.61%  '
$ws.Range('E38').Value = '  -6.39%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.73'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  -1.00%  '
$ws.Range('E40').Value = '  -3.05%  '
$ws.Range('E41').Value = '  +4.83%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '69.42'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  -0.91%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.228'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  -0.11%  '
$ws.Range('E44').Value = '  -0.16%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '93.39'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  -3.56%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.28'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  +1.58%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '1.741.19'
$ws.Range('E47').Value = '  +9.12%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '114.08'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  +1.81%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '80.35'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  -1.49%  '
$ws.Range('E50').Value = '  -1.63%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.05'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  -3.50%  '
